$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: exposure date corrected from 30/12/20 to 21/12/20
$ws.Range("C3").Value = "21/12/20 1pm-1:30pm"

# Row 29 (Moorabin / Grape and Grain Liquor Cellars entry) is removed entirely;
# deleting the whole row shifts rows 30-31 up to become rows 29-30 and updates
# the sheet dimension from A1:D31 to A1:D30 automatically.
$ws.Rows(29).Delete()

# Row 27 (still "Melbourne" in column A) now gets new site/period/notes describing
# a Pakenham line train exposure instead of the Melbourne Central Lion Hotel entry.
$ws.Range("B27").Value = "Metro Train - Pakenham line"
$ws.Range("C27").Value = "01/01/21 3:00am-4:00am"
$ws.Range("D27").Value = "Case caught train from Flinders St to Westall Station"
